$d = $word.ActiveDocument

# The published edit narrows to a single word swap inside the "Summary:" line:
#   "An independent and self-motivated student ..." -> "... self-motivated employee ..."
# The canonical XML shows the run that used to read
#   "An independent and self-motivated student with proven ability and experienced in"
# split into three runs, with only the middle one ("employee") differing from the
# original ("student"). Re-create that by locating just the word "student" and
# replacing it in place (a direct Range.Text assignment on the whole original run
# would just swap the whole sentence into one run, losing the 3-way split that a
# real "select word -> type replacement" edit produces in Word).

$r = $d.Content
$found = $r.Find.Execute("student", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the word 'student' to replace."
}

# Toggling a character-format property on this narrowed Range forces Word to split
# the run at the Range boundaries (exactly what happens when you select a word and
# retype it). Flipping it back off afterwards restores the original (unformatted)
# look while keeping the new run boundaries.
$r.Bold = 1
$r.Text = "employee"
$r.Bold = 0
